$wb = $excel.ActiveWorkbook

# Add new worksheet after "Disposition_master" (2nd sheet), before "Updation_of_Disposition"
$afterSheet = $wb.Worksheets.Item("Disposition_master")
$newSheet = $wb.Worksheets.Add($null, $afterSheet)
$newSheet.Name = "Agency_Account_Allocation"

$newSheet.Range("A1").Value = "TestScenario"
$newSheet.Range("B1").Value = "Run"
$newSheet.Range("C1").Value = "DPD"
$newSheet.Range("A2").Value = "Agency_Account_Allocation"
$newSheet.Range("B2").Value = "Yes"
$newSheet.Range("C2").Value = 250

$newSheet.Activate()
